$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.405.72"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "1.845.26"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'264.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.5200"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("D8").Value = "'0.3269"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("D9").Value = "'0.06806"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'18.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("D11").Value = "'0.7757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'0.07781"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.831.64"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "'87.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'13.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007987"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "26.421.63"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "2.064.02"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "'4.635"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'9.567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("D24").Value = "'5.990"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").Value = "'2.183"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.73%  "
$ws.Range("D27").Value = "'1.658"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "'17.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").Value = "'112.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'4.177"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").Value = "'4.136"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").Value = "'0.08744"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").Value = "'0.04826"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").Value = "'0.7198"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "'2.855"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "'3.095"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").Value = "'0.01779"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("D40").Value = "'0.4862"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.96%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("D42").Value = "'111.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").Value = "'6.060"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'7.711"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'0.4161"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").Value = "'9.097"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "'0.1235"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.56%  "
$ws.Range("D50").Value = "'34.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.93%  "
$ws.Range("D51").Value = "'0.8884"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.47%  "
